$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.424.35"
$ws.Range("E2").Value = "'  +1.01%  "

$ws.Range("D3").Value = "'1.903.25"
$ws.Range("E3").Value = "'  +2.39%  "

$ws.Range("E4").Value = "'  +0.28%  "

$ws.Range("D5").Value = "'246.10"
$ws.Range("E5").Value = "'  +3.80%  "

$ws.Range("D6").Value = "'0.638"
$ws.Range("E6").Value = "'  +2.58%  "

$ws.Range("E7").Value = "'  +0.26%  "

$ws.Range("E8").Value = "'  -1.71%  "

$ws.Range("E9").Value = "'  +3.76%  "

$ws.Range("D10").Value = "'0.0705"
$ws.Range("E10").Value = "'  +1.27%  "

$ws.Range("E11").Value = "'  +0.91%  "

$ws.Range("D12").Value = "'2.181.97"
$ws.Range("E12").Value = "'  +2.68%  "

$ws.Range("D13").Value = "'12.37"
$ws.Range("E13").Value = "'  +8.30%  "

$ws.Range("D14").Value = "'0.698"
$ws.Range("E14").Value = "'  +2.70%  "

$ws.Range("D15").Value = "'1.897.24"
$ws.Range("E15").Value = "'  +1.84%  "

$ws.Range("E16").Value = "'  +2.54%  "

$ws.Range("D17").Value = "'35.486.16"
$ws.Range("E17").Value = "'  +1.29%  "

$ws.Range("D18").Value = "'71.90"
$ws.Range("E18").Value = "'  +2.27%  "

$ws.Range("D19").Value = "'0.0₃0828"
$ws.Range("E19").Value = "'  +3.96%  "

$ws.Range("D20").Value = "'243.20"

$ws.Range("D21").Value = "'12.67"
$ws.Range("E21").Value = "'  +4.06%  "

$ws.Range("D22").Value = "'4.82"
$ws.Range("E22").Value = "'  +1.43%  "

$ws.Range("E23").Value = "'  +0.17%  "

$ws.Range("E24").Value = "'  +0.50%  "

$ws.Range("B25").Value = "'PancakeSwap"
$ws.Range("C25").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "'  +14.34%  "

$ws.Range("B26").Value = "'Monero"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'171.72"
$ws.Range("E26").Value = "'  +0.18%  "

$ws.Range("D27").Value = "'8.54"
$ws.Range("E27").Value = "'  +7.74%  "

$ws.Range("D28").Value = "'18.00"
$ws.Range("E28").Value = "'  +1.83%  "

$ws.Range("E29").Value = "'  +0.87%  "

$ws.Range("D30").Value = "'0.976"
$ws.Range("E30").Value = "'  +24.92%  "

$ws.Range("D31").Value = "'0.0569"
$ws.Range("E31").Value = "'  +2.14%  "

$ws.Range("E32").Value = "'  +2.79%  "

$ws.Range("E33").Value = "'  +0.29%  "

$ws.Range("D34").Value = "'4.16"
$ws.Range("E34").Value = "'  +3.95%  "

$ws.Range("E35").Value = "'  +8.70%  "

$ws.Range("D36").Value = "'2.02"
$ws.Range("E36").Value = "'  -0.15%  "

$ws.Range("E37").Value = "'  +4.75%  "

$ws.Range("D38").Value = "'1.10"
$ws.Range("E38").Value = "'  +1.95%  "

$ws.Range("B39").Value = "'Kaspa"
$ws.Range("C39").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.0639"
$ws.Range("E39").Value = "'  +16.20%  "

$ws.Range("D40").Value = "'91.92"
$ws.Range("E40").Value = "'  +0.25%  "

$ws.Range("B41").Value = "'VeChain"
$ws.Range("C41").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0204"
$ws.Range("E41").Value = "'  +0.82%  "

$ws.Range("D42").Value = "'15.59"
$ws.Range("E42").Value = "'  +4.51%  "

$ws.Range("D43").Value = "'1.346.66"
$ws.Range("E43").Value = "'  -0.45%  "

$ws.Range("D44").Value = "'50.31"
$ws.Range("E44").Value = "'  +44.81%  "

$ws.Range("E45").Value = "'  +2.39%  "

$ws.Range("D46").Value = "'12.93"
$ws.Range("E46").Value = "'  +1.48%  "

$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "'  +0.00%  "

$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = "'  -0.15%  "

$ws.Range("D49").Value = "'6.61"
$ws.Range("E49").Value = "'  +4.40%  "

$ws.Range("D50").Value = "'2.092.74"
$ws.Range("E50").Value = "'  +2.36%  "

$ws.Range("D51").Value = "'0.0691"
$ws.Range("E51").Value = "'  +1.53%  "
